$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sku"
$ws.Range("B2").Value = "name"
$ws.Range("C2").Value = "quantity"
$ws.Range("D2").Value = "cost_per"
$ws.Range("E2").Value = "total_cost"

$ws.Range("A3").Value = "sku"
$ws.Range("B3").Value = "name"
$ws.Range("C3").Value = "quantity"
$ws.Range("D3").Value = "cost_per"
$ws.Range("E3").Value = "total_cost"

$ws.Range("A4").Value = "sku"
$ws.Range("B4").Value = "name"
$ws.Range("C4").Value = "quantity"
$ws.Range("D4").Value = "cost_per"
$ws.Range("E4").Value = "total_cost"

$ws.Range("A5").Value = "sku"
$ws.Range("B5").Value = "name"
$ws.Range("C5").Value = "quantity"
$ws.Range("D5").Value = "cost_per"
$ws.Range("E5").Value = "total_cost"
